$d = $word.ActiveDocument

$replacements = @(
    @("61×89=", "73×32="),
    @("36×53=", "66×44="),
    @("32×16=", "84×57="),
    @("80×72=", "80×61="),
    @("79×23=", "52×71="),
    @("76×79=", "35×17="),
    @("62×71=", "91×16="),
    @("31×31=", "42×11="),
    @("71×92=", "14×99="),
    @("81×62=", "55×12="),
    @("97×49=", "77×11="),
    @("82×16=", "63×62="),
    @("15×94=", "64×61="),
    @("25×21=", "65×68="),
    @("34×54=", "66×57="),
    @("33×62=", "94×26="),
    @("19×47=", "55×23="),
    @("13×82=", "31×37="),
    @("78×54=", "78×11="),
    @("94×75=", "88×88="),
    @("78×51=", "19×37="),
    @("93×90=", "91×26="),
    @("66×61=", "85×51="),
    @("13×76=", "89×70="),
    @("13×25=", "40×61=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
